$d = $word.ActiveDocument

# --- Replacement 1 ---
$old1 = "SMID, P. CNC Programming Handbook, Industrial Press, 2007.GROOVER, M.; ZIMEMERS, E. Computer Aided Design and Manufacturing, Prentice-Hall, 1984.STENERSON, J.; CURRAN, K. Computer Numerical Control: Operation and Programming, Prentice Hall, 2006.SIMON, W. Numerical Control of Machine Tools, Edward Arnold, 1973.MILNER, D.; VASILOV, V.: Computer Aided Engineering for Manufacture. Kogan Page, 1986.CHUA, C. K.; LEONG, K. F. Rapid Prototyping: Principles and Applications, World Scientific Publishing, 2010. MESSLER, R. W. Joining of Materials and Structures, Butterworth-Heinemann, 2004.KIMINAMI, C. S.; CASTRO, W. B.; OLIVEIRA, M. F. Introdução aos processos de Fabricação de Produtos Metálicos, Blucher, 2013.MEYERS, M.A. AND CHAWLA, K.K.; Mechanical Behavior of Materials, Prentice-Hall, Upper Saddle River-NJ (EUA), 1999.GIESECKE, F. E. Comunicação Gráfica Moderna, Editora Bookman, 2002.CRUZ, M. D. Catia V5r20 - Modelagem, Montagem e Detalhamento, ERICA, 2010.FISCHER, U; GOMERINGER, R; HEINZLER, M; ET AL. Manual de Tecnologia Metal Mecânica, Blucher, 2011.JACK, H. Projeto, Planejamento e Gestão de Produtos: Uma abordagem para engenharia, Campus-Elsevier, 2014.SWIFT, K.G.; BOOKER, P.D. Seleção de processos de manufatura, Campus-Elsevier, 2014."
$new1 = "SMID, P. CNC Programming Handbook, Industrial Press, 2007.^lGROOVER, M.; ZIMEMERS, E. Computer Aided Design and Manufacturing, Prentice-Hall, 1984.^lSTENERSON, J.; CURRAN, K. Computer Numerical Control: Operation and Programming, Prentice Hall, 2006.^lSIMON, W. Numerical Control of Machine Tools, Edward Arnold, 1973.^lMILNER, D.; VASILOV, V.: Computer Aided Engineering for Manufacture. Kogan Page, 1986.^lCHUA, C. K.; LEONG, K. F. Rapid Prototyping: Principles and Applications, World Scientific Publishing, 2010. MESSLER, R. W. Joining of Materials and Structures, Butterworth-Heinemann, 2004.^lKIMINAMI, C. S.; CASTRO, W. B.; OLIVEIRA, M. F. Introdução aos processos de Fabricação de Produtos Metálicos, Blucher, 2013.^lMEYERS, M.A. AND CHAWLA, K.K.; Mechanical Behavior of Materials, Prentice-Hall, Upper Saddle River-NJ (EUA), 1999.^lGIESECKE, F. E. Comunicação Gráfica Moderna, Editora Bookman, 2002.^lCRUZ, M. D. Catia V5r20 - Modelagem, Montagem e Detalhamento, ERICA, 2010.^lFISCHER, U; GOMERINGER, R; HEINZLER, M; ET AL. Manual de Tecnologia Metal Mecânica, Blucher, 2011.^lJACK, H. Projeto, Planejamento e Gestão de Produtos: Uma abordagem para engenharia, Campus-Elsevier, 2014.^lSWIFT, K.G.; BOOKER, P.D. Seleção de processos de manufatura, Campus-Elsevier, 2014."
$found1 = $d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)
Write-Host "Replacement 1 found: $found1"


# --- Replacement 2 ---
$old2 = "O projeto extensionista proposto visa capacitar a comunidade local em práticas seguras e sustentáveis de fabricação, com foco em temas relacionados à Engenharia Física., incluindo possibilidade de fabricação de protótipos. Através de oficinas de construção utilizando ferramentas manuais e Equipamentos de Proteção Individual (EPIs), os alunos do curso serão protagonistas no desenvolvimento e execução das atividades, enquanto a comunidade será contemplada com conhecimentos importantes para a segurança no trabalho e a conscientização ambiental.Através de oficinas práticas, os participantes serão capacitados a adotar medidas de segurança no trabalho e a incorporar princípios de sustentabilidade em suas atividades de fabricação, contribuindo para um ambiente de trabalho mais seguro e para a promoção do desenvolvimento sustentável na comunidade."
$new2 = "O projeto extensionista proposto visa capacitar a comunidade local em práticas seguras e sustentáveis de fabricação, com foco em temas relacionados à Engenharia Física., incluindo possibilidade de fabricação de protótipos. Através de oficinas de construção utilizando ferramentas manuais e Equipamentos de Proteção Individual (EPIs), os alunos do curso serão protagonistas no desenvolvimento e execução das atividades, enquanto a comunidade será contemplada com conhecimentos importantes para a segurança no trabalho e a conscientização ambiental.^lAtravés de oficinas práticas, os participantes serão capacitados a adotar medidas de segurança no trabalho e a incorporar princípios de sustentabilidade em suas atividades de fabricação, contribuindo para um ambiente de trabalho mais seguro e para a promoção do desenvolvimento sustentável na comunidade."
$found2 = $d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)
Write-Host "Replacement 2 found: $found2"


# --- Replacement 3 ---
$old3 = "Os alunos do curso de Engenharia Física serão os protagonistas no desenvolvimento e execução das atividades do projeto. As etapas incluirão:•Planejamento e Preparação: Os alunos irão pesquisar e selecionar os conteúdos a serem abordados nas oficinas, além de preparar o material didático necessário.•Divulgação e Inscrição: Divulgação do projeto na comunidade através de panfletos, redes sociais e parcerias com instituições locais, parcerias com escolas, associações comunitárias e empresas locais. As inscrições serão abertas para os interessados em participar das oficinas.•Realização das Oficinas: As oficinas serão realizadas em um local adequado, com equipamentos de segurança e materiais necessários fornecidos. Os alunos serão responsáveis por ministrar as aulas práticas, abordando os temas propostos.•Avaliação dos Participantes: Ao final de cada oficina, os participantes serão avaliados quanto ao conhecimento adquirido, habilidades práticas desenvolvidas e percepção sobre a importância da segurança e sustentabilidade na fabricação.•Feedback e Melhoria Contínua: Os alunos irão coletar feedback dos participantes e utilizar essas informações para aprimorar as próximas edições das oficinas."
$new3 = "Os alunos do curso de Engenharia Física serão os protagonistas no desenvolvimento e execução das atividades do projeto. As etapas incluirão:^l•Planejamento e Preparação: Os alunos irão pesquisar e selecionar os conteúdos a serem abordados nas oficinas, além de preparar o material didático necessário.^l•Divulgação e Inscrição: Divulgação do projeto na comunidade através de panfletos, redes sociais e parcerias com instituições locais, parcerias com escolas, associações comunitárias e empresas locais. As inscrições serão abertas para os interessados em participar das oficinas.^l•Realização das Oficinas: As oficinas serão realizadas em um local adequado, com equipamentos de segurança e materiais necessários fornecidos. Os alunos serão responsáveis por ministrar as aulas práticas, abordando os temas propostos.^l•Avaliação dos Participantes: Ao final de cada oficina, os participantes serão avaliados quanto ao conhecimento adquirido, habilidades práticas desenvolvidas e percepção sobre a importância da segurança e sustentabilidade na fabricação.^l•Feedback e Melhoria Contínua: Os alunos irão coletar feedback dos participantes e utilizar essas informações para aprimorar as próximas edições das oficinas."
$found3 = $d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2)
Write-Host "Replacement 3 found: $found3"


# --- Replacement 4 ---
$old4 = "•Questionários de satisfação ao final de cada oficina.•Acompanhamento do número de participantes que adotam práticas seguras e sustentáveis em seus trabalhos após a participação nas oficinas.•Questionários de satisfação ao final de cada oficina.•Acompanhamento do número de participantes que adotam práticas seguras e sustentáveis em seus trabalhos após a participação nas oficinas.•Feedback verbal durante as interações com os participantes (sinalizar aqui como o grupo social externo à Universidade poderá avaliar a atividade realizada conjuntamente com os estudantes, durante sua realização e ao final)"
$new4 = "•Questionários de satisfação ao final de cada oficina.^l•Acompanhamento do número de participantes que adotam práticas seguras e sustentáveis em seus trabalhos após a participação nas oficinas.^l•Questionários de satisfação ao final de cada oficina.^l•Acompanhamento do número de participantes que adotam práticas seguras e sustentáveis em seus trabalhos após a participação nas oficinas.^l•Feedback verbal durante as interações com os participantes (sinalizar aqui como o grupo social externo à Universidade poderá avaliar a atividade realizada conjuntamente com os estudantes, durante sua realização e ao final)"
$found4 = $d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $new4, 2)
Write-Host "Replacement 4 found: $found4"
